$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 3 ("Mira Road East, Thane" / "Gyms" row).
# Excel will automatically shift all subsequent rows up by one.
$ws.Rows.Item(3).Delete()

# Re-sequence the S.No column (A) so it stays a contiguous 1..8 list
# after the row removal.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Reflect the last cell selection recorded after the edit.
$ws.Range("B14").Select()
